$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new sale line item ---

# A7: counter goes from 0 to 1
$ws.Range("A7").Value = 1

# C7:G7 (merged) - item name -> General/Text cell, plain string assignment
# becomes a shared-string text cell; also switch number format to Text (@)
# to match the workbook's updated formatting for this column (whole merged
# range gets reformatted, same as the underlying shared cell style).
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "TORSERETIC 20MG 30 TABS."

# H7:K7 (merged) - quantity/ratio text "0:2"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "0:2"

# L7:M7 (merged) - "1" stored as text, but keep its original (numeric-looking)
# number format untouched, same as the source workbook.
$origFmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $origFmtL7

# N7:O7 (merged) - price text "123.00"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "123.00"

# P7 - "123.0000" stored as text, keep the original numeric format untouched.
$origFmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "123.0000"
$ws.Range("P7").NumberFormat = $origFmtP7

# Q7 - "1:0" text
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Row 8: total row ---
$ws.Range("P8").Value = 123

# --- Row 9: footer (timestamp refreshed, page + credit text re-saved) ---
$ws.Range("A9").Value = "Wednesday, 23 July, 2025 9:20 AM"
$ws.Range("G9").Value = "1/1"
$ws.Range("K9").Value = "developed by : Abdelaziz Talaat"
